$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("table")

# --- Row 7: P7 "Tem and Rep Only" -> "Yes" (now using the regular/unhighlighted style) ---
$ws.Range("P7").Value2 = "Yes"
$ws.Range("P7").Font.Bold = $false

# New note in S7 (Notes column) for row 7
$ws.Range("S7").Value2 = "Synapsida Replicate 5 marginal_rates error"

# --- Rows 12-15: fill in the "Name in Files" column (A) which had been left blank ---
$ws.Range("A12").Value2 = "mcmc_predictors/B_covar_rjmcmc"
$ws.Range("A13").Value2 = "mcmc_predictors/B_covar_rjmcmc"
$ws.Range("A14").Value2 = "mcmc_predictors/B_covar_rjmcmc"
$ws.Range("A15").Value2 = "mcmc_predictors/B_covar_rjmcmc"

# Rows 14 & 15: mark "Submitted?" (O) as concluded/done -> "Yes" (unbolded, keep fill)
$ws.Range("O14").Value2 = "Yes"
$ws.Range("O14").Font.Bold = $false
$ws.Range("O15").Value2 = "Yes"
$ws.Range("O15").Font.Bold = $false

# --- Row 29: update / shorten the footnote text and its row height ---
$ws.Range("B29").Value2 = "*All Models use a Time Variable Poisson Process (TPP) to model Preservation"
$ws.Rows.Item(29).RowHeight = 51.45

# --- Row 22: mark this run as restored (append the "*" footnote reference) ---
$ws.Range("A22").Value2 = "mcmc_no_predictors/A_bdnn_update & A_bdnn*"

# --- Restore selection on the sheet to D25 ---
$ws.Activate()
$ws.Range("D25").Select()
